# Atualização - RDD 13 e Copas
# Adiciona a coluna "Rodada 13" (N) com as pontuações da rodada, mantendo
# as linhas/times existentes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cabeçalho da nova rodada (copia o formato do cabeçalho anterior - negrito,
# centralizado e com borda - e então define o texto)
$ws.Range("M1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "Rodada 13"

# Pontuações da Rodada 13 por time (linha)
$ws.Range("N3").Value  = 107.64990234375
$ws.Range("N4").Value  = 96.509765625
$ws.Range("N5").Value  = 70.2099609375
$ws.Range("N6").Value  = 118.740234375
$ws.Range("N8").Value  = 110.41015625
$ws.Range("N9").Value  = 117.509765625
$ws.Range("N10").Value = 128.0400390625
$ws.Range("N12").Value = 93.509765625
$ws.Range("N13").Value = 111.91015625
$ws.Range("N14").Value = 117.10986328125
$ws.Range("N15").Value = 133.0498046875
$ws.Range("N19").Value = 113.509765625
$ws.Range("N21").Value = 79.41015625
$ws.Range("N22").Value = 120.8701171875
$ws.Range("N26").Value = 119.41015625
$ws.Range("N32").Value = 124.31005859375
$ws.Range("N33").Value = 130.2099609375
